$wb = $excel.ActiveWorkbook

# Move "Grafici per le slides" so that it sits right before "Sheet1"
# (i.e. swap the order of these two sheets; "Mapping campi" stays first).
$sheet1 = $wb.Worksheets.Item("Sheet1")
$grafici = $wb.Worksheets.Item("Grafici per le slides")
$grafici.Move($sheet1)

# Re-fetch "Sheet1" by name (worksheet handles are position based in this
# runtime, so grab a fresh reference after the reorder) and rename it.
$sheet1 = $wb.Worksheets.Item("Sheet1")
$sheet1.Name = "Suddivisione boxes in the slide"

# Make the renamed sheet the active tab, matching the edited workbook
# (3rd tab / activeTab index 2 selected).
$target = $wb.Worksheets.Item("Suddivisione boxes in the slide")
$target.Activate()
